$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "Layer by Layer" boolean rows first (rows 3,5,7 each
# duplicate the probability-distribution already present in the row above),
# so that only one row per probability distribution remains.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(3).Delete()

# Drop the "Experiment" (col A) and "Layer by Layer" (col C) columns, shifting
# "Probability distribution" into column A and "CIFAR-10/MNIST test accuracy"
# into columns B/C.
$ws.Columns.Item(3).Delete()
$ws.Columns.Item(1).Delete()

# Reset the selection back to the top-left cell (matches the saved view).
$ws.Range("A1").Select() | Out-Null
